$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 391, shifting existing rows 391-461 down to 392-462.
$ws.Rows.Item(391).Insert()

# Populate the newly inserted row 391 with the new data record.
$ws.Cells.Item(391, 1).Value = 1
$ws.Cells.Item(391, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(391, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(391, 4).Value = 45275
$ws.Cells.Item(391, 5).Value = 15
$ws.Cells.Item(391, 6).Value = "Fruta"
$ws.Cells.Item(391, 7).Value = 100108
$ws.Cells.Item(391, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(391, 9).Value = 100108006
$ws.Cells.Item(391, 10).Value = "Plátano"
$ws.Cells.Item(391, 11).Value = "Sin especificar"
$ws.Cells.Item(391, 12).Value = "Pintón"
$ws.Cells.Item(391, 13).Value = 200
$ws.Cells.Item(391, 14).Value = 20000
$ws.Cells.Item(391, 15).Value = 22000
$ws.Cells.Item(391, 16).Value = 21000
$ws.Cells.Item(391, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(391, 18).Value = "Ecuador"
$ws.Cells.Item(391, 19).Value = 1050
$ws.Cells.Item(391, 20).Value = 20
